# Apply crime-index YTD data update for 2023-11-12
# For each affected worksheet, update the listed cells (by A1-style reference)
# with their new values.

$wb = $excel.ActiveWorkbook

$updates = [ordered]@{
    'Citywide Totals' = @{ "I2"=110; "D3"=126; "J3"=210; "J4"=2; "B9"=352; "C9"=452; "D9"=383; "E9"=427; "H9"=420; "I9"=474; "J9"=385; "B10"=1237; "C10"=1456; "D10"=1677; "E10"=1982; "F10"=1970; "G10"=855; "I10"=795; "J10"=667; "B11"=1714; "C11"=2061; "D11"=2282; "E11"=2624; "F11"=2666; "G11"=1495; "H11"=1232; "I11"=1590; "J11"=1400 }
    'By Neighborhood' = @{ "D5"=16; "F8"=133; "G10"=12; "D18"=6; "B27"=19; "J27"=20; "E28"=85; "F28"=113; "J28"=57; "C32"=94; "E32"=135; "I32"=92; "D36"=74; "I41"=15; "I42"=10; "F43"=16; "E50"=45; "B53"=239; "C53"=345; "D53"=552; "E53"=658; "G53"=236; "J53"=220; "H61"=22; "J61"=4; "F63"=12; "D65"=62; "E65"=44; "J67"=4; "F70"=100; "I72"=14; "D74"=92; "E76"=91; "J77"=69; "F78"=50; "E92"=30; "I92"=34; "D96"=31; "B97"=28; "H98"=7; "B99"=1714; "C99"=2061; "D99"=2282; "E99"=2624; "F99"=2666; "G99"=1495; "H99"=1232; "I99"=1590; "J99"=1400 }
    'Rogers Park' = @{ "E7"=73; "E8"=91 }
    'Edgewater' = @{ "B4"=3; "J5"=9; "B6"=19; "J6"=20 }
    'Roseland' = @{ "J8"=18; "J10"=69 }
    'Austin' = @{ "F7"=91; "F8"=133 }
    'Garfield Park' = @{ "C7"=33; "I7"=30; "C8"=56; "E8"=76; "C9"=94; "E9"=135; "I9"=92 }
    'Grand Crossing' = @{ "D3"=8; "D9"=74 }
    'Loop' = @{ "J4"=2; "E7"=61; "B8"=193; "C8"=293; "D8"=485; "E8"=581; "G8"=156; "B9"=239; "C9"=345; "D9"=552; "E9"=658; "G9"=236; "J9"=220 }
    'Armour Square' = @{ "D6"=9; "D7"=16 }
    'Old Town' = @{ "F7"=77; "F8"=100 }
    'Little Italy, UIC' = @{ "E5"=18; "E7"=45 }
    'North Lawndale' = @{ "D6"=20; "D7"=41; "E7"=33; "D8"=62; "E8"=44 }
    'Humboldt Park' = @{ "I2"=1; "I6"=15 }
    'Rush & Division' = @{ "F5"=40; "F6"=50 }
    'Englewood' = @{ "E8"=55; "F8"=62; "J8"=28; "E9"=85; "F9"=113; "J9"=57 }
    'River North' = @{ "D6"=81; "D7"=92 }
    'West Loop' = @{ "I7"=6; "E8"=25; "E9"=30; "I9"=34 }
    'Norwood Park' = @{ "J5"=2; "J6"=4 }
    'Woodlawn' = @{ "B5"=10; "B7"=28 }
    'Calumet Heights' = @{ "D5"=4; "D6"=6 }
    'Wicker Park' = @{ "D6"=28; "D7"=31 }
    'New City' = @{ "F5"=7; "F6"=12 }
    'Printers Row' = @{ "I4"=9; "I6"=14 }
    'Wrigleyville' = @{ "H5"=1; "H7"=7 }
    'Hyde Park' = @{ "I5"=5; "I6"=10 }
    'Avondale' = @{ "G6"=8; "G7"=12 }
    'Irving Park' = @{ "F6"=13; "F7"=16 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
